$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.617.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +4.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.825.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +5.38%  '
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '340.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.19%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3830'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.34%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3537'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.75'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.241'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07750'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("E13").Value = '  +10.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.631'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.76%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.822.80'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.211'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("E17").Value = '  +5.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06727'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.29'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9996'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.554'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.16'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.29%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.607.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.477'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.684'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.28%  '
$ws.Range("E27").Value = '  +13.82%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.488'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.66%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '152.93'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.028.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.49'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.360'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.084'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08820'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.697'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.09%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.653'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.17%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7059'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.169'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.50%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2270'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.07%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06527'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.65%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02412'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.298'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.83'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +5.21%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6639'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +10.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.927'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.30%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.191'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.32%  '
$ws.Range("E49").Value = '  +4.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07313'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.69'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.84%  '
